$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")
$ws.Activate()

# Row 8: new expense entry, matching the date-cell formatting used by rows above (A2:A7)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = 66434

$ws.Range("B8").Value = "chiya patti 1pau, biscuit 1 packet, advance for asset ring"
$ws.Range("C8").Formula = "=110+80+10000"

# Update current selection on the "all" sheet to D7
$ws.Range("D7").Select()
